$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = 8.199999999999999
$ws.Range("AC2").Value = 8.4
$ws.Range("AJ2").Value = 24
$ws.Range("L2").Value = 1.43
$ws.Range("P2").Value = 1.84
$ws.Range("Q2").Value = 2.06
$ws.Range("S2").Value = 3.7
$ws.Range("T2").Value = 1.9
$ws.Range("U2").Value = 1.99
$ws.Range("AA3").Value = 20
$ws.Range("AB3").Value = 16.5
$ws.Range("AC3").Value = 8.4
$ws.Range("AD3").Value = 10.5
$ws.Range("AE3").Value = 28
$ws.Range("AF3").Value = 50
$ws.Range("AG3").Value = 30
$ws.Range("AH3").Value = 980
$ws.Range("AJ3").Value = 190
$ws.Range("AL3").Value = 120
$ws.Range("AM3").Value = 180
$ws.Range("AO3").Value = 14.5
$ws.Range("G3").Value = 5.6
$ws.Range("H3").Value = 1.79
$ws.Range("N3").Value = 3.3
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 1.78
$ws.Range("Q3").Value = 2.18
$ws.Range("R3").Value = 1.3
$ws.Range("U3").Value = 1.89
$ws.Range("X3").Value = 12.5
$ws.Range("Y3").Value = 7.4
$ws.Range("Z3").Value = 10.5
$ws.Range("AC4").Value = 8.6
$ws.Range("AO4").Value = 18
$ws.Range("F4").Value = 5.5
$ws.Range("G4").Value = 5.7
$ws.Range("I4").Value = 1.85
$ws.Range("K4").Value = 3.6
$ws.Range("P4").Value = 1.62
$ws.Range("U4").Value = 1.74
$ws.Range("AA5").Value = 980
$ws.Range("AB5").Value = 18.5
$ws.Range("AD5").Value = 10.5
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 980
$ws.Range("AH5").Value = 32
$ws.Range("AI5").Value = 44
$ws.Range("AO5").Value = 13.5
$ws.Range("P5").Value = 1.88
$ws.Range("T5").Value = 1.95
$ws.Range("U5").Value = 1.98
$ws.Range("Z5").Value = 11
$ws.Range("F6").Value = 3.2
$ws.Range("H6").Value = 2.32
$ws.Range("I6").Value = 2.5
$ws.Range("AN7").Value = 70
$ws.Range("F7").Value = 3.3
$ws.Range("N7").Value = 3.75
$ws.Range("R7").Value = 1.37
$ws.Range("AH8").Value = 21
$ws.Range("AI8").Value = 1000
$ws.Range("AJ8").Value = 17.5
$ws.Range("AK8").Value = 20
$ws.Range("AO8").Value = 1000
$ws.Range("F8").Value = 1.71
$ws.Range("G8").Value = 1.73
$ws.Range("I8").Value = 5.9
$ws.Range("J8").Value = 4
$ws.Range("K8").Value = 4.6
$ws.Range("N8").Value = 4.2
$ws.Range("P8").Value = 2.14
$ws.Range("Q8").Value = 1.72
$ws.Range("S8").Value = 2.52
$ws.Range("U8").Value = 2.06
$ws.Range("Z8").Value = 55
